$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original plain-text representation
# (e.g. "1.00", "0.610", "0.0000263") instead of Excel auto-coercing these
# numeric-looking strings into native numbers (which would drop trailing
# zeros / reformat them). Apply a text format across the data rows once.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '64.789.22'
$ws.Range('E2').Value = '  -3.81%  '
$ws.Range('D3').Value = '3.341.36'
$ws.Range('E3').Value = '  -4.58%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '182.21'
$ws.Range('E5').Value = '  -8.90%  '
$ws.Range('D6').Value = '534.93'
$ws.Range('E6').Value = '  -2.90%  '
$ws.Range('D7').Value = '0.610'
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('D8').Value = '3.337.35'
$ws.Range('E8').Value = '  -4.39%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  -5.26%  '
$ws.Range('D11').Value = '59.01'
$ws.Range('E11').Value = '  -7.23%  '
$ws.Range('D12').Value = '0.136'
$ws.Range('E12').Value = '  -5.34%  '
$ws.Range('D13').Value = '0.0000263'
$ws.Range('E13').Value = '  -2.43%  '
$ws.Range('E14').Value = '  -6.37%  '
$ws.Range('D15').Value = '3.869.15'
$ws.Range('E15').Value = '  -4.82%  '
$ws.Range('D16').Value = '3.340.38'
$ws.Range('E16').Value = '  -4.61%  '
$ws.Range('E17').Value = '  -4.36%  '
$ws.Range('D18').Value = '64.700.62'
$ws.Range('E18').Value = '  -3.60%  '
$ws.Range('D19').Value = '17.68'
$ws.Range('E19').Value = '  -3.50%  '
$ws.Range('D20').Value = '11.25'
$ws.Range('E20').Value = '  -4.33%  '
$ws.Range('D21').Value = '0.971'
$ws.Range('E21').Value = '  -4.92%  '
$ws.Range('D22').Value = '378.28'
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('D23').Value = '3.85'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').Value = '11.30'
$ws.Range('E24').Value = '  -6.53%  '
$ws.Range('D25').Value = '81.37'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('E26').Value = '  +2.13%  '
$ws.Range('E27').Value = '  -1.10%  '
$ws.Range('D28').Value = '2.71'
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('D29').Value = '11.58'
$ws.Range('E29').Value = '  -5.07%  '
$ws.Range('D30').Value = '8.47'
$ws.Range('E30').Value = '  -3.74%  '
$ws.Range('D31').Value = '29.27'
$ws.Range('E31').Value = '  -5.42%  '
$ws.Range('D32').Value = '660.30'
$ws.Range('E33').Value = '  -2.82%  '
$ws.Range('D34').Value = '11.39'
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('D36').Value = '59.83'
$ws.Range('E36').Value = '  -6.19%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '0.396'
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').Value = '37.22'
$ws.Range('E39').Value = '  -3.89%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0715'
$ws.Range('E40').Value = '  +6.49%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = '0.126'
$ws.Range('E42').Value = '  -3.45%  '
$ws.Range('D43').Value = '2.947.49'
$ws.Range('E43').Value = '  -3.86%  '
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('E45').Value = '  -7.73%  '
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('E47').Value = '  -3.66%  '
$ws.Range('E48').Value = '  +6.03%  '
$ws.Range('D49').Value = '2.81'
$ws.Range('E49').Value = '  +7.28%  '
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('E51').Value = '  -5.08%  '
